$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = 3292
$ws.Range("E2").Value = 688
$ws.Range("F2").Value = 688
$ws.Range("G2").Value = 617
$ws.Range("H2").Value = 571
$ws.Range("I2").Value = 571
$ws.Range("K2").Value = 3709
$ws.Range("L2").Value = 734
$ws.Range("M2").Value = 2975
$ws.Range("N2").Value = 2975
$ws.Range("P2").Value = 86
$ws.Range("Q2").Value = 468
$ws.Range("R2").Value = -690
$ws.Range("S2").Value = 450
$ws.Range("T2").Value = 108
$ws.Range("U2").Value = 361
$ws.Range("V2").Value = 205
$ws.Range("W2").Value = 20.91
$ws.Range("X2").Value = 17.36
$ws.Range("Y2").Value = 26.35
$ws.Range("Z2").Value = 17.9
$ws.Range("AA2").Value = 24.68
$ws.Range("AB2").Value = 3681.8
$ws.Range("AC2").Value = 3914
$ws.Range("AD2").Value = 6.76
$ws.Range("AE2").Value = 18268
$ws.Range("AF2").Value = 1.45
$ws.Range("AG2").Value = 700
$ws.Range("AH2").Value = 2.65
$ws.Range("AI2").Value = 14.5
$ws.Range("AJ2").Value = 17285715
$ws.Range("J2").ClearContents()
$ws.Range("O2").ClearContents()

# Row 3
$ws.Range("D3").Value = 2971
$ws.Range("E3").Value = 364
$ws.Range("F3").Value = 364
$ws.Range("G3").Value = 258
$ws.Range("H3").Value = 214
$ws.Range("I3").Value = 214
$ws.Range("J3").Value = 0
$ws.Range("K3").Value = 3920
$ws.Range("L3").Value = 663
$ws.Range("M3").Value = 3257
$ws.Range("N3").Value = 3255
$ws.Range("O3").Value = 2
$ws.Range("P3").Value = 86
$ws.Range("Q3").Value = -198
$ws.Range("R3").Value = 638
$ws.Range("S3").Value = -204
$ws.Range("T3").Value = 174
$ws.Range("U3").Value = -373
$ws.Range("V3").Value = 150
$ws.Range("W3").Value = 12.26
$ws.Range("X3").Value = 7.2
$ws.Range("Y3").Value = 6.86
$ws.Range("Z3").Value = 5.61
$ws.Range("AA3").Value = 20.36
$ws.Range("AB3").Value = 3834.35
$ws.Range("AC3").Value = 1236
$ws.Range("AD3").Value = 20.15
$ws.Range("AE3").Value = 19936
$ws.Range("AF3").Value = 1.25
$ws.Range("AG3").Value = 300
$ws.Range("AH3").Value = 1.2
$ws.Range("AI3").Value = 22.93
$ws.Range("AJ3").Value = 17285715

# Row 4
$ws.Range("D4").Value = 3112
$ws.Range("E4").Value = -237
$ws.Range("F4").Value = -237
$ws.Range("G4").Value = -168
$ws.Range("H4").Value = -172
$ws.Range("I4").Value = -173
$ws.Range("J4").Value = 0
$ws.Range("K4").Value = 3857
$ws.Range("L4").Value = 871
$ws.Range("M4").Value = 2986
$ws.Range("N4").Value = 2984
$ws.Range("O4").Value = 2
$ws.Range("P4").Value = 86
$ws.Range("Q4").Value = 488
$ws.Range("R4").Value = -330
$ws.Range("S4").Value = -125
$ws.Range("T4").Value = 160
$ws.Range("U4").Value = 329
$ws.Range("V4").Value = 169
$ws.Range("W4").Value = -7.61
$ws.Range("X4").Value = -5.54
$ws.Range("Y4").Value = -5.54
$ws.Range("Z4").Value = -4.44
$ws.Range("AA4").Value = 29.15
$ws.Range("AB4").Value = 3576.42
$ws.Range("AC4").Value = -999
$ws.Range("AD4").Value = -17.57
$ws.Range("AE4").Value = 18749
$ws.Range("AF4").Value = 0.9399999999999999
$ws.Range("AG4").Value = 200
$ws.Range("AH4").Value = 1.14
$ws.Range("AI4").Value = -18.43
$ws.Range("AJ4").Value = 17285715

# Row 5
$ws.Range("D5").Value = 3122
$ws.Range("E5").Value = 354
$ws.Range("F5").Value = 354
$ws.Range("G5").Value = 419
$ws.Range("H5").Value = 390
$ws.Range("I5").Value = 389
$ws.Range("J5").Value = 0
$ws.Range("K5").Value = 3879
$ws.Range("L5").Value = 861
$ws.Range("M5").Value = 3019
$ws.Range("N5").Value = 3017
$ws.Range("O5").Value = 2
$ws.Range("P5").Value = 86
$ws.Range("Q5").Value = -31
$ws.Range("R5").Value = -139
$ws.Range("S5").Value = 95
$ws.Range("T5").Value = 367
$ws.Range("U5").Value = -398
$ws.Range("V5").Value = 300
$ws.Range("W5").Value = 11.34
$ws.Range("X5").Value = 12.47
$ws.Range("Y5").Value = 12.98
$ws.Range("Z5").Value = 10.07
$ws.Range("AA5").Value = 28.51
$ws.Range("AB5").Value = 3986.87
$ws.Range("AC5").Value = 2253
$ws.Range("AD5").Value = 11.63
$ws.Range("AE5").Value = 19109
$ws.Range("AF5").Value = 1.37
$ws.Range("AG5").Value = 350
$ws.Range("AH5").Value = 1.34
$ws.Range("AI5").Value = 14.19
$ws.Range("AJ5").Value = 17285715

# Row 6
$ws.Range("D6").Value = 5022
$ws.Range("E6").Value = 327
$ws.Range("F6").Value = 327
$ws.Range("G6").Value = 68
$ws.Range("H6").Value = 71
$ws.Range("I6").Value = 59
$ws.Range("K6").Value = 6064
$ws.Range("L6").Value = 2879
$ws.Range("M6").Value = 3185
$ws.Range("N6").Value = 3021
$ws.Range("P6").Value = 86
$ws.Range("Q6").Value = -205
$ws.Range("R6").Value = -617
$ws.Range("S6").Value = 793
$ws.Range("T6").Value = 264
$ws.Range("U6").Value = -469
$ws.Range("V6").Value = 1777
$ws.Range("W6").Value = 6.52
$ws.Range("X6").Value = 1.41
$ws.Range("Y6").Value = 1.96
$ws.Range("Z6").Value = 1.42
$ws.Range("AA6").Value = 90.38
$ws.Range("AB6").Value = 4081.14
$ws.Range("AC6").Value = 343
$ws.Range("AD6").Value = 80.53
$ws.Range("AE6").Value = 19520
$ws.Range("AF6").Value = 1.41
$ws.Range("AG6").Value = 400
$ws.Range("AH6").Value = 1.45
$ws.Range("AI6").Value = 104.49
$ws.Range("AJ6").Value = 17285715

# Row 7
$ws.Range("D7").Value = 7802
$ws.Range("E7").Value = 606
$ws.Range("G7").Value = 546
$ws.Range("H7").Value = 452
$ws.Range("I7").Value = 445
$ws.Range("K7").Value = 7224
$ws.Range("L7").Value = 3576
$ws.Range("M7").Value = 3648
$ws.Range("N7").Value = 3469
$ws.Range("P7").Value = 86
$ws.Range("Q7").Value = 592
$ws.Range("R7").Value = -377
$ws.Range("S7").Value = -329
$ws.Range("T7").Value = 401
$ws.Range("U7").Value = -49
$ws.Range("W7").Value = 7.76
$ws.Range("X7").Value = 5.8
$ws.Range("Y7").Value = 13.71
$ws.Range("Z7").Value = 6.81
$ws.Range("AA7").Value = 98.03
$ws.Range("AC7").Value = 2574
$ws.Range("AD7").Value = 14.29
$ws.Range("AE7").Value = 22396
$ws.Range("AF7").Value = 1.64
$ws.Range("AG7").Value = 400
$ws.Range("AH7").Value = 1.09
$ws.Range("AI7").Value = 15.54

# Row 8
$ws.Range("D8").Value = 9396
$ws.Range("E8").Value = 838
$ws.Range("G8").Value = 738
$ws.Range("H8").Value = 632
$ws.Range("I8").Value = 636
$ws.Range("K8").Value = 8145
$ws.Range("L8").Value = 3946
$ws.Range("M8").Value = 4199
$ws.Range("N8").Value = 4020
$ws.Range("P8").Value = 86
$ws.Range("Q8").Value = 645
$ws.Range("R8").Value = -471
$ws.Range("S8").Value = -69
$ws.Range("T8").Value = 480
$ws.Range("U8").Value = 175
$ws.Range("W8").Value = 8.92
$ws.Range("X8").Value = 6.73
$ws.Range("Y8").Value = 16.97
$ws.Range("Z8").Value = 8.220000000000001
$ws.Range("AA8").Value = 93.97
$ws.Range("AC8").Value = 3676
$ws.Range("AD8").Value = 10.01
$ws.Range("AE8").Value = 25953
$ws.Range("AF8").Value = 1.42
$ws.Range("AG8").Value = 400
$ws.Range("AH8").Value = 1.09
$ws.Range("AI8").Value = 10.88

# Row 9
$ws.Range("D9").Value = 10322
$ws.Range("E9").Value = 978
$ws.Range("G9").Value = 891
$ws.Range("H9").Value = 751
$ws.Range("I9").Value = 754
$ws.Range("K9").Value = 9010
$ws.Range("L9").Value = 4142
$ws.Range("M9").Value = 4868
$ws.Range("N9").Value = 4689
$ws.Range("P9").Value = 86
$ws.Range("Q9").Value = 744
$ws.Range("R9").Value = -520
$ws.Range("S9").Value = -69
$ws.Range("T9").Value = 527
$ws.Range("U9").Value = 225
$ws.Range("W9").Value = 9.48
$ws.Range("X9").Value = 7.27
$ws.Range("Y9").Value = 17.33
$ws.Range("Z9").Value = 8.76
$ws.Range("AA9").Value = 85.09
$ws.Range("AC9").Value = 4365
$ws.Range("AD9").Value = 8.43
$ws.Range("AE9").Value = 30272
$ws.Range("AF9").Value = 1.22
$ws.Range("AG9").Value = 400
$ws.Range("AH9").Value = 1.09
$ws.Range("AI9").Value = 9.16
